# Update Name of Algo
# Apply updated imputation result values produced by the RandomForest run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value  = 12.73779999999999
$ws.Range("E6").Value  = 11.9327
$ws.Range("D7").Value  = -7.450499999999998
$ws.Range("E7").Value  = 12.6415
$ws.Range("B8").Value  = 5.369799999999998
$ws.Range("E8").Value  = 14.01979999999999
$ws.Range("E9").Value  = 9.741499999999988
$ws.Range("B10").Value = 8.726000000000001
$ws.Range("E10").Value = 11.3371
$ws.Range("B12").Value = 6.252100000000002
$ws.Range("E12").Value = 12.95399999999999
$ws.Range("C13").Value = -12.87999999999999
$ws.Range("B18").Value = 4.998100000000004
$ws.Range("D20").Value = -8.068300000000001
$ws.Range("B25").Value = 5.745999999999993
